# Auto-generated script to apply numeric corrections to the Ultros_Profits leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1435.3158
$ws.Range("J107").Value = 401.22223
$ws.Range("L107").Value = 401.22223
$ws.Range("N107").Value = -4241.22223

$ws.Range("H113").Value = 9947.315000000001
$ws.Range("I113").Value = 10000
$ws.Range("J113").Value = 9928.5
$ws.Range("K113").Value = 10000
$ws.Range("L113").Value = 9928.5
$ws.Range("M113").Value = -6746
$ws.Range("N113").Value = -16436.5

$ws.Range("H131").Value = 13607.5
$ws.Range("I131").Value = 14765.833
$ws.Range("J131").Value = 12449.167
$ws.Range("K131").Value = 44297.499
$ws.Range("L131").Value = 37347.501
$ws.Range("M131").Value = -39257.499
$ws.Range("N131").Value = -47427.501

$ws.Range("H138").Value = 3079
$ws.Range("I138").Value = 1804.1428
$ws.Range("K138").Value = 5412.428400000001
$ws.Range("M138").Value = -272.4284000000007

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 14827.808
$ws.Range("I2").Value = 17846.3
$ws.Range("J2").Value = 4766.1665
$ws.Range("K2").Value = 17846.3
$ws.Range("L2").Value = 4766.1665
$ws.Range("M2").Value = -17733.3
$ws.Range("N2").Value = -4992.1665

$ws.Range("H45").Value = 3907.9
$ws.Range("I45").Value = 2513.25
$ws.Range("K45").Value = 2513.25
$ws.Range("M45").Value = -2136.25

$ws.Range("H61").Value = 1839
$ws.Range("I61").Value = 1839
$ws.Range("K61").Value = 1839
$ws.Range("M61").Value = -1627

$ws.Range("H74").Value = 1762.25
$ws.Range("I74").Value = 1474.7084
$ws.Range("J74").Value = 2624.875
$ws.Range("K74").Value = 1474.7084
$ws.Range("L74").Value = 2624.875
$ws.Range("M74").Value = -600.7084
$ws.Range("N74").Value = -4372.875

$ws.Range("H77").Value = 1762.25
$ws.Range("I77").Value = 1474.7084
$ws.Range("J77").Value = 2624.875
$ws.Range("K77").Value = 7373.541999999999
$ws.Range("L77").Value = 13124.375
$ws.Range("M77").Value = -3005.541999999999
$ws.Range("N77").Value = -21860.375

$ws.Range("H110").Value = 5078
$ws.Range("I110").Value = 4498.9414
$ws.Range("K110").Value = 4498.9414
$ws.Range("M110").Value = -2453.9414

$ws.Range("H116").Value = 14827.808
$ws.Range("I116").Value = 17846.3
$ws.Range("J116").Value = 4766.1665
$ws.Range("K116").Value = 17846.3
$ws.Range("L116").Value = 4766.1665
$ws.Range("M116").Value = -15552.3
$ws.Range("N116").Value = -9354.166499999999

$ws.Range("H122").Value = 6375
$ws.Range("I122").Value = 6000
$ws.Range("K122").Value = 18000
$ws.Range("M122").Value = -15550

$ws.Range("H136").Value = 1839
$ws.Range("I136").Value = 1839
$ws.Range("K136").Value = 5517
$ws.Range("M136").Value = -2967

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 14827.808
$ws.Range("I3").Value = 17846.3
$ws.Range("J3").Value = 4766.1665
$ws.Range("K3").Value = 17846.3
$ws.Range("L3").Value = 4766.1665
$ws.Range("M3").Value = -17732.3
$ws.Range("N3").Value = -4994.1665

$ws.Range("H20").Value = 2428.5293
$ws.Range("I20").Value = 1884
$ws.Range("J20").Value = 2809.7
$ws.Range("K20").Value = 1884
$ws.Range("L20").Value = 2809.7
$ws.Range("M20").Value = -1637
$ws.Range("N20").Value = -3303.7

$ws.Range("H22").Value = 33533.668
$ws.Range("I22").Value = 33533.668
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 33533.668
$ws.Range("N22").Value = 0
$ws.Range("M22").Value = -33360.668
$ws.Range("L22").ClearContents()

$ws.Range("H25").Value = 400
$ws.Range("I25").Value = 400
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 400
$ws.Range("L25").Value = 0
$ws.Range("N25").Value = -165
$ws.Range("M25").ClearContents()

$ws.Range("H94").Value = 3060.9
$ws.Range("I94").Value = 2682.3845
$ws.Range("J94").Value = 3763.8572
$ws.Range("K94").Value = 2682.3845
$ws.Range("L94").Value = 3763.8572
$ws.Range("M94").Value = -2231.3845
$ws.Range("N94").Value = -4665.8572

$ws.Range("H105").Value = 1737.6666
$ws.Range("I105").Value = 1555.4166
$ws.Range("K105").Value = 1555.4166
$ws.Range("M105").Value = 191.5834

$ws.Range("H134").Value = 2532.261
$ws.Range("I134").Value = 2138
$ws.Range("J134").Value = 2962.3635
$ws.Range("K134").Value = 6414
$ws.Range("L134").Value = 8887.0905
$ws.Range("M134").Value = -3879
$ws.Range("N134").Value = -13957.0905

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1270.8334
$ws.Range("I5").Value = 1333.3334
$ws.Range("K5").Value = 1333.3334
$ws.Range("M5").Value = -1221.3334

$ws.Range("H16").Value = 1341.4286
$ws.Range("I16").Value = 1195.4286
$ws.Range("J16").Value = 1487.4286
$ws.Range("K16").Value = 1195.4286
$ws.Range("L16").Value = 1487.4286
$ws.Range("M16").Value = -908.4286
$ws.Range("N16").Value = -2061.4286

$ws.Range("H41").Value = 3055.5557
$ws.Range("I41").Value = 3055.5557
$ws.Range("K41").Value = 3055.5557
$ws.Range("M41").Value = -2627.5557

$ws.Range("H47").Value = 40000
$ws.Range("J47").Value = 40000
$ws.Range("L47").Value = 40000
$ws.Range("N47").Value = -41132

$ws.Range("H60").Value = 294.6111
$ws.Range("I60").Value = 294.6111
$ws.Range("K60").Value = 294.6111
$ws.Range("M60").Value = 216.3889

$ws.Range("H113").Value = 1341.4286
$ws.Range("I113").Value = 1195.4286
$ws.Range("J113").Value = 1487.4286
$ws.Range("K113").Value = 1195.4286
$ws.Range("L113").Value = 1487.4286
$ws.Range("M113").Value = 974.5714
$ws.Range("N113").Value = -5827.4286

$ws.Range("H117").Value = 71599.8
$ws.Range("J117").Value = 71599.8
$ws.Range("L117").Value = 71599.8
$ws.Range("N117").Value = -80777.8

$ws.Range("H122").Value = 4390.75
$ws.Range("I122").Value = 4133
$ws.Range("J122").Value = 4545.4
$ws.Range("K122").Value = 12399
$ws.Range("L122").Value = 13636.2
$ws.Range("M122").Value = -9949
$ws.Range("N122").Value = -18536.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 851.625
$ws.Range("I8").Value = 851.625
$ws.Range("K8").Value = 2554.875
$ws.Range("M8").Value = -2415.875

$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()

$ws.Range("H17").Value = 150
$ws.Range("I17").Value = 150
$ws.Range("K17").Value = 450
$ws.Range("M17").Value = -281

$ws.Range("H34").Value = 2901.875
$ws.Range("J34").Value = 2816.4285
$ws.Range("L34").Value = 8449.2855
$ws.Range("N34").Value = -8617.2855

$ws.Range("H39").Value = 2488.7273
$ws.Range("J39").Value = 3196.5715
$ws.Range("L39").Value = 9589.7145
$ws.Range("N39").Value = -10177.7145

$ws.Range("H40").Value = 108.666664
$ws.Range("J40").Value = 136
$ws.Range("L40").Value = 544
$ws.Range("N40").Value = -682

$ws.Range("H55").Value = 2850.9666
$ws.Range("J55").Value = 4252.5
$ws.Range("L55").Value = 12757.5
$ws.Range("N55").Value = -13111.5

$ws.Range("H92").Value = 622.2222
$ws.Range("I92").Value = 694.5
$ws.Range("J92").Value = 601.5714
$ws.Range("K92").Value = 2083.5
$ws.Range("L92").Value = 1804.7142
$ws.Range("M92").Value = -835.5
$ws.Range("N92").Value = -4300.7142

$ws.Range("H122").Value = 6228
$ws.Range("I122").Value = 913.5
$ws.Range("J122").Value = 7999.5
$ws.Range("K122").Value = 8221.5
$ws.Range("L122").Value = 71995.5
$ws.Range("M122").Value = -5771.5
$ws.Range("N122").Value = -76895.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 5050000
$ws.Range("I7").Value = 10000000
$ws.Range("K7").Value = 10000000
$ws.Range("M7").Value = -9999888

$ws.Range("H8").Value = 5050000
$ws.Range("I8").Value = 10000000
$ws.Range("K8").Value = 10000000
$ws.Range("M8").Value = -9999861

$ws.Range("H122").Value = 9471.143
$ws.Range("I122").Value = 6166.6665
$ws.Range("K122").Value = 18499.9995
$ws.Range("M122").Value = -16049.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 839.8
$ws.Range("I19").Value = 675
$ws.Range("J19").Value = 1499
$ws.Range("K19").Value = 675
$ws.Range("L19").Value = 1499
$ws.Range("M19").Value = -505
$ws.Range("N19").Value = -1839

$ws.Range("H93").Value = 2589.5833
$ws.Range("I93").Value = 1980.5555
$ws.Range("K93").Value = 1980.5555
$ws.Range("M93").Value = -732.5554999999999

$ws.Range("H100").Value = 52683.914
$ws.Range("I100").Value = 68895.88
$ws.Range("K100").Value = 68895.88
$ws.Range("M100").Value = -68354.88

$ws.Range("H122").Value = 6839.5
$ws.Range("I122").Value = 6452.6665
$ws.Range("K122").Value = 19357.9995
$ws.Range("M122").Value = -16907.9995

$ws.Range("H136").Value = 4249.0527
$ws.Range("I136").Value = 3793.818
$ws.Range("J136").Value = 4875
$ws.Range("K136").Value = 11381.454
$ws.Range("L136").Value = 14625
$ws.Range("M136").Value = -8831.454000000002
$ws.Range("N136").Value = -19725

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 6500
$ws.Range("J22").Value = 6500
$ws.Range("L22").Value = 6500
$ws.Range("N22").Value = -7086

$ws.Range("H24").Value = 23800
$ws.Range("J24").Value = 23800
$ws.Range("L24").Value = 23800
$ws.Range("N24").Value = -24260

$ws.Range("H100").Value = 1122
$ws.Range("I100").Value = 1324.8334
$ws.Range("J100").Value = 716.3333
$ws.Range("K100").Value = 2649.6668
$ws.Range("L100").Value = 1432.6666
$ws.Range("M100").Value = -2108.6668
$ws.Range("N100").Value = -2514.6666

$ws.Range("H122").Value = 1306.7
$ws.Range("I122").Value = 1286.3334
$ws.Range("K122").Value = 3859.0002
$ws.Range("M122").Value = -1409.0002

